# "3 cases check in"
# - VerifyCSVForNewVersion: clear the H column "Pass" results for rows 5-58
# - VerifyEventAPI: fill in the H column results for rows 5-51 ("Pass", except
#   row 51 which is "Fail")
# - BatchDecisionOutputValidations: remove the stray F7/F8/F9 result cells

$wb = $excel.ActiveWorkbook

# --- VerifyCSVForNewVersion: clear H5:H58 ---
$wsNew = $wb.Worksheets.Item("VerifyCSVForNewVersion")
for ($r = 5; $r -le 58; $r++) {
    $wsNew.Cells.Item($r, 8).Value = ""
}

# --- VerifyEventAPI: set H5:H51 to "Pass", with H51 = "Fail" ---
$wsEvt = $wb.Worksheets.Item("VerifyEventAPI")
for ($r = 5; $r -le 50; $r++) {
    $wsEvt.Cells.Item($r, 8).Value = "Pass"
}
$wsEvt.Cells.Item(51, 8).Value = "Fail"

# --- BatchDecisionOutputValidations: drop F7, F8, F9 ---
$wsBatch = $wb.Worksheets.Item("BatchDecisionOutputValidations")
$wsBatch.Cells.Item(7, 6).Value = ""
$wsBatch.Cells.Item(8, 6).Value = ""
$wsBatch.Cells.Item(9, 6).Value = ""
